$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update status of "5. Bank transaction" (row 21) and
# "6. Petty Cash, sama dengan The Wave" (row 22) to "Done"
$ws.Range("B21").Value = "Done"
$ws.Range("B22").Value = "Done"

# Scroll the view back to the top and move the selection to E17,
# as captured when the sheet was last saved.
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("E17").Select()
